$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 28 to make room for the new Step13/14/15
# rows (25-27) and keep the Methods table (myMethod1..6) pushed down by 2.
$ws.Rows("28:29").Insert()

# New "Steps" entries: explicit casts to/from alias datatypes and array casts.
$ws.Range("C25").Value = "Step13"
$ws.Range("D25").Value = "'= (Datatype4) 1"
$ws.Range("C26").Value = "Step14"
$ws.Range("D26").Value = "'= (Datatype4[]) `$Step13"
$ws.Range("C27").Value = "Step15"
$ws.Range("D27").Value = "'= (Datatype5[]) `$Step13"

# Restore the exact original column widths (engine's ColumnWidth is applied
# in points, so compensate for its internal padding to land back on 21/14).
$ws.Columns("C").ColumnWidth = 20.166666666666668
$ws.Columns("E").ColumnWidth = 13.166666666666666

# Update sheet view / selection to match the saved state after editing.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I26").Select()
